# Scheduled-runner update: refresh profit-tracking figures across the
# Sheets workbook (ALC/ARM/CRP/CUL/GSM/LTW/WVR) with newly computed
# pricing data. Pure value updates, no formulas / formatting touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 787.05884
$ws.Range("I28").Value = 825.7273
$ws.Range("J28").Value = 716.1667
$ws.Range("K28").Value = 825.7273
$ws.Range("L28").Value = 716.1667
$ws.Range("M28").Value = -340.7273
$ws.Range("N28").Value = -1686.1667

$ws.Range("H51").Value = 3282.9524
$ws.Range("J51").Value = 4485.5
$ws.Range("L51").Value = 4485.5
$ws.Range("N51").Value = -5453.5

$ws.Range("H76").Value = 3229.5334
$ws.Range("I76").Value = 3116.1775
$ws.Range("K76").Value = 3116.1775
$ws.Range("M76").Value = -2801.1775

$ws.Range("H79").Value = 3229.5334
$ws.Range("I79").Value = 3116.1775
$ws.Range("K79").Value = 3116.1775
$ws.Range("M79").Value = -2024.1775

$ws.Range("H116").Value = 23549.5
$ws.Range("I116").Value = 6698.4287
$ws.Range("K116").Value = 6698.4287
$ws.Range("M116").Value = -3256.4287

$ws.Range("H132").Value = 8665.933999999999
$ws.Range("I132").Value = 8713.5
$ws.Range("K132").Value = 26140.5
$ws.Range("M132").Value = -23610.5

$ws.Range("H137").Value = 3834.681
$ws.Range("I137").Value = 1574.5834
$ws.Range("J137").Value = 11231.363
$ws.Range("K137").Value = 4723.7502
$ws.Range("L137").Value = 33694.089
$ws.Range("M137").Value = -2173.7502
$ws.Range("N137").Value = -38794.089

$ws.Range("H141").Value = 6641
$ws.Range("I141").Value = 3217.5454
$ws.Range("K141").Value = 9652.636200000001
$ws.Range("M141").Value = -4472.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7747.2856
$ws.Range("I31").Value = 6538.5
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 6538.5
$ws.Range("L31").Value = 15000
$ws.Range("M31").Value = -6244.5
$ws.Range("N31").Value = -15588

$ws.Range("H32").Value = 159272.33
$ws.Range("I32").Value = 179801.77
$ws.Range("K32").Value = 179801.77
$ws.Range("M32").Value = -179514.77

$ws.Range("H61").Value = 3341.9697
$ws.Range("I61").Value = 3735.6365
$ws.Range("J61").Value = 2554.6365
$ws.Range("K61").Value = 3735.6365
$ws.Range("L61").Value = 2554.6365
$ws.Range("M61").Value = -3523.6365
$ws.Range("N61").Value = -2978.6365

$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248

$ws.Range("H63").Value = 3500
$ws.Range("I63").Value = 3500
$ws.Range("K63").Value = 3500
$ws.Range("M63").Value = -2814

$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240

$ws.Range("H66").Value = 3500
$ws.Range("I66").Value = 3500
$ws.Range("K66").Value = 17500
$ws.Range("M66").Value = -14068

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 13234.889
$ws.Range("I74").Value = 1398.909
$ws.Range("K74").Value = 1398.909
$ws.Range("M74").Value = -524.9090000000001

$ws.Range("H77").Value = 13234.889
$ws.Range("I77").Value = 1398.909
$ws.Range("K77").Value = 6994.545
$ws.Range("M77").Value = -2626.545

$ws.Range("H102").Value = 4333
$ws.Range("I102").Value = 4303.737
$ws.Range("J102").Value = 4611
$ws.Range("K102").Value = 4303.737
$ws.Range("L102").Value = 4611
$ws.Range("M102").Value = -2681.737
$ws.Range("N102").Value = -7855

$ws.Range("H110").Value = 885.1667
$ws.Range("I110").Value = 686.4400000000001
$ws.Range("K110").Value = 686.4400000000001
$ws.Range("M110").Value = 1358.56

$ws.Range("H122").Value = 2066.7856
$ws.Range("I122").Value = 1975.7693
$ws.Range("K122").Value = 5927.3079
$ws.Range("M122").Value = -3477.3079

$ws.Range("H132").Value = 927129.7
$ws.Range("I132").Value = 1042751.44
$ws.Range("K132").Value = 3128254.32
$ws.Range("M132").Value = -3125724.32

$ws.Range("H136").Value = 3341.9697
$ws.Range("I136").Value = 3735.6365
$ws.Range("J136").Value = 2554.6365
$ws.Range("K136").Value = 11206.9095
$ws.Range("L136").Value = 7663.9095
$ws.Range("M136").Value = -8656.9095
$ws.Range("N136").Value = -12763.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 48150.43
$ws.Range("I16").Value = 516.8461
$ws.Range("K16").Value = 516.8461
$ws.Range("M16").Value = -229.8461

$ws.Range("H31").Value = 2503.1924
$ws.Range("I31").Value = 2495.75
$ws.Range("J31").Value = 2509.5715
$ws.Range("K31").Value = 2495.75
$ws.Range("L31").Value = 2509.5715
$ws.Range("M31").Value = -2200.75
$ws.Range("N31").Value = -3099.5715

$ws.Range("H34").Value = 2503.1924
$ws.Range("I34").Value = 2495.75
$ws.Range("J34").Value = 2509.5715
$ws.Range("K34").Value = 2495.75
$ws.Range("L34").Value = 2509.5715
$ws.Range("M34").Value = -2293.75
$ws.Range("N34").Value = -2913.5715

$ws.Range("H99").Value = 8270358.5
$ws.Range("I99").Value = 32941.855
$ws.Range("K99").Value = 32941.855
$ws.Range("M99").Value = -31443.855

$ws.Range("H107").Value = 683
$ws.Range("I107").Value = 640
$ws.Range("K107").Value = 640
$ws.Range("M107").Value = 1280

$ws.Range("H113").Value = 48150.43
$ws.Range("I113").Value = 516.8461
$ws.Range("K113").Value = 516.8461
$ws.Range("M113").Value = 1653.1539

$ws.Range("H126").Value = 8270358.5
$ws.Range("I126").Value = 32941.855
$ws.Range("K126").Value = 98825.565
$ws.Range("M126").Value = -96355.565

$ws.Range("H134").Value = 1688.45
$ws.Range("I134").Value = 1542.7778
$ws.Range("K134").Value = 4628.3334
$ws.Range("M134").Value = -2093.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 55.25
$ws.Range("I2").Value = 44.2
$ws.Range("K2").Value = 265.2
$ws.Range("M2").Value = -152.2

$ws.Range("H101").Value = 8346083
$ws.Range("J101").Value = 8346083
$ws.Range("L101").Value = 25038249
$ws.Range("N101").Value = -25043117

$ws.Range("H113").Value = 1512.5
$ws.Range("I113").Value = 692.3333
$ws.Range("K113").Value = 2076.9999
$ws.Range("M113").Value = 93.0001000000002

$ws.Range("H121").Value = 2868.389
$ws.Range("I121").Value = 512
$ws.Range("J121").Value = 6571.2856
$ws.Range("K121").Value = 1536
$ws.Range("L121").Value = 19713.8568
$ws.Range("M121").Value = -226
$ws.Range("N121").Value = -22333.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 14666.667
$ws.Range("J52").Value = 14666.667
$ws.Range("L52").Value = 14666.667
$ws.Range("N52").Value = -15184.667

$ws.Range("H80").Value = 1916.6666

$ws.Range("H83").Value = 1916.6666

$ws.Range("H97").Value = 1277.8914
$ws.Range("I97").Value = 1239.5641
$ws.Range("K97").Value = 1239.5641
$ws.Range("M97").Value = -743.5641000000001

$ws.Range("H113").Value = 2703.2
$ws.Range("I113").Value = 1841.8096
$ws.Range("K113").Value = 1841.8096
$ws.Range("M113").Value = 328.1904

$ws.Range("H117").Value = 86650
$ws.Range("J117").Value = 86650
$ws.Range("L117").Value = 86650
$ws.Range("N117").Value = -93534

$ws.Range("H122").Value = 2552
$ws.Range("I122").Value = 2041.6957
$ws.Range("J122").Value = 3856.111
$ws.Range("K122").Value = 6125.0871
$ws.Range("L122").Value = 11568.333
$ws.Range("M122").Value = -3675.0871
$ws.Range("N122").Value = -16468.333

$ws.Range("H132").Value = 16682.848
$ws.Range("I132").Value = 20880.916
$ws.Range("K132").Value = 62642.74800000001
$ws.Range("M132").Value = -60112.74800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2911.6538
$ws.Range("I100").Value = 2657.2856
$ws.Range("J100").Value = 3980
$ws.Range("K100").Value = 2657.2856
$ws.Range("L100").Value = 3980
$ws.Range("M100").Value = -2116.2856
$ws.Range("N100").Value = -5062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 8000
$ws.Range("J19").Value = 8000
$ws.Range("L19").Value = 8000
$ws.Range("N19").Value = -8348

$ws.Range("H55").Value = 13333.333
$ws.Range("I55").Value = 10000
$ws.Range("J55").Value = 15000
$ws.Range("K55").Value = 10000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -9723
$ws.Range("N55").Value = -15554

$ws.Range("H126").Value = 2794.6667
$ws.Range("I126").Value = 2763.077
$ws.Range("K126").Value = 8289.231
$ws.Range("M126").Value = -5819.231
